# Update the "想去人数" (F column) figures for several entries on the
# "展览" and "全部类型" worksheets, matching the new data snapshot.

$wb = $excel.ActiveWorkbook

# Row -> new F-column value for both affected sheets.
$updates = @{
    11 = 4249
    13 = 279
    17 = 60
    18 = 3001
    19 = 58
    20 = 429
    31 = 412
    32 = 1673
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}

$wb.Save()
